# Figuras atualizadas manualmente antes de criar definir a atualização automática.
#
# The sheet holds one row per (Região, Ano) pair for the "Gastos públicos
# com segurança" series, grouped by region (Sergipe, Nordeste, Brasil) and
# ordered by year within each group. A new year (2024) is being added to
# every region, so one new row is inserted right after the last existing
# row of each region's block (pushing the following blocks down), and is
# filled in with that region's 2024 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sergipe: insert the 2024 row after the existing 2015-2023 rows (2-10) ---
$ws.Rows.Item(11).Insert()
$ws.Cells.Item(11, 1).Value = "Sergipe"
$ws.Cells.Item(11, 2).Value = "Gastos públicos com segurança"
$ws.Cells.Item(11, 3).NumberFormat = "@"
$ws.Cells.Item(11, 3).Value = "01/01/2024"
$ws.Cells.Item(11, 4).Value = 768.9122732866518
$ws.Cells.Item(11, 5).Value = 9

# --- Nordeste: insert the 2024 row after the existing 2015-2023 rows (now 12-20) ---
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).Value = "Nordeste"
$ws.Cells.Item(21, 2).Value = "Gastos públicos com segurança"
$ws.Cells.Item(21, 3).NumberFormat = "@"
$ws.Cells.Item(21, 3).Value = "01/01/2024"
$ws.Cells.Item(21, 4).Value = 538.2672412097861
$ws.Cells.Item(21, 5).Value = ""

# --- Brasil: insert the 2024 row after the existing 2015-2023 rows (now 22-30) ---
$ws.Rows.Item(31).Insert()
$ws.Cells.Item(31, 1).Value = "Brasil"
$ws.Cells.Item(31, 2).Value = "Gastos públicos com segurança"
$ws.Cells.Item(31, 3).NumberFormat = "@"
$ws.Cells.Item(31, 3).Value = "01/01/2024"
$ws.Cells.Item(31, 4).Value = 711.5263293354096
$ws.Cells.Item(31, 5).Value = ""
